# Brazil.xlsx monthly M2 update
# - corrects B419/D419, B420/D420, B421/D421 (revised historical values)
# - appends three new monthly rows (422-424) for 2023-07, 2023-08, 2023-09

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrections to existing rows ---
$ws.Cells.Item(419, 2).Value = 5245350000000
$ws.Cells.Item(419, 4).Value = 1051425192431.045

$ws.Cells.Item(420, 2).Value = 5322265000000
$ws.Cells.Item(420, 4).Value = 1052788107765.953

$ws.Cells.Item(421, 2).Value = 5421638000000
$ws.Cells.Item(421, 4).Value = 1131984131955.319

# --- New rows 422-424 ---
$ws.Cells.Item(422, 1).Value = 45108
$ws.Cells.Item(422, 2).Value = 5501072000000
$ws.Cells.Item(422, 3).Value = 0.2116894938504202
$ws.Cells.Item(422, 4).Value = 1164519147314.719

$ws.Cells.Item(423, 1).Value = 45139
$ws.Cells.Item(423, 2).Value = 5591097000000
$ws.Cells.Item(423, 3).Value = 0.2018733850129199
$ws.Cells.Item(423, 4).Value = 1128693677325.582

$ws.Cells.Item(424, 1).Value = 45170
$ws.Cells.Item(424, 2).Value = 5656835000000
$ws.Cells.Item(424, 3).Value = 0.1987281399046105
$ws.Cells.Item(424, 4).Value = 1124172297297.297

# Column A carries the date number format / border style (s="2") used
# throughout the sheet; replicate it onto the new rows by copying the
# format from the row directly above.
$ws.Range("A421").Copy()
$ws.Range("A422:A424").PasteSpecial(-4122)
$excel.CutCopyMode = $false
